# Update "lines_states" worksheet: add line7/line8 rows, shift extr rows down,
# and update the C/D/E values as per the new contingency computations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write one data row (A..E) given row index.
function Set-DataRow($r, $aVal, $bName, $cVal, $dVal, $eVal) {
    $ws.Cells.Item($r, 1).Value2 = $aVal
    $ws.Cells.Item($r, 2).Value2 = $bName
    $ws.Cells.Item($r, 3).Value2 = $cVal
    $ws.Cells.Item($r, 4).Value2 = $dVal
    $ws.Cells.Item($r, 5).Value2 = $eVal
}

# Copy the style of an already-styled index cell (A14, which has the bold/border style)
# onto the two brand-new rows (16 and 17) before they hold any data, so no new
# style entries get created in styles.xml.
$ws.Range("A14").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 2-7 are unchanged.

# Row 8: was extr1 -> now line7
Set-DataRow 8 6 "line7" 14 11 $true

# Row 9: was extr2 -> now line8
Set-DataRow 9 7 "line8" 16 9 $true

# Row 10: was extr3 -> now extr1
Set-DataRow 10 8 "extr1" 5 12 $true

# Row 11: was extr4 -> now extr2
Set-DataRow 11 9 "extr2" 5 9 $true

# Row 12: was extr5 -> now extr3
Set-DataRow 12 10 "extr3" 10 11 $false

# Row 13: was extr6 -> now extr4
Set-DataRow 13 11 "extr4" 7 8 $true

# Row 14: was extr7 -> now extr5
Set-DataRow 14 12 "extr5" 9 11 $false

# Row 15: was extr8 -> now extr6
Set-DataRow 15 13 "extr6" 7 11 $false

# Row 16 (new): extr7
Set-DataRow 16 14 "extr7" 5 7 $true

# Row 17 (new): extr8
Set-DataRow 17 15 "extr8" 8 5 $false

Write-Host "Done updating lines_states sheet"
